$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2: D0.04039067246480155 -> D0.9894230751189594
$ws.Range("A2").Value = "D0.9894230751189594"

# A3: D0.22709280606726223 -> D0.6966197969744036
$ws.Range("A3").Value = "D0.6966197969744036"

# A4: D0.3037386565467328 -> D0.6600436572746871
$ws.Range("A4").Value = "D0.6600436572746871"

# D2, D3, D4: Sat, 10 Dec 2022 21:22:08 +0530 -> Sat, 24 Dec 2022 00:45:23 -0800
$ws.Range("D2").Value = "Sat, 24 Dec 2022 00:45:23 -0800"
$ws.Range("D3").Value = "Sat, 24 Dec 2022 00:45:23 -0800"
$ws.Range("D4").Value = "Sat, 24 Dec 2022 00:45:23 -0800"
